$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.104.46"
$ws.Range("E2").Value = "  -3.31%  "
$ws.Range("D3").Value = "1.601.80"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3782"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08158"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.599"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001259"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.419"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.33%  "
$ws.Range("D17").Value = "1.593.33"
$ws.Range("E17").Value = "  -3.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06851"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.585"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5542"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.17%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.68%  "
$ws.Range("D25").Value = "23.105.14"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.340"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.715"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.17%  "
$ws.Range("E28").Value = "  -4.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.283"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.396"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.867"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -12.87%  "
$ws.Range("D34").Value = "1.778.84"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9632"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07706"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.292"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.25%  "
$ws.Range("E38").Value = "  -6.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2559"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08903"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.370"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7092"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6626"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.316"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.994"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07939"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.40%  "
